# "Added test for constant vol surface."
#
# The existing "vol_surface" sheet (Tenors / Quotes / blank helper column)
# gets duplicated into a brand-new sheet "constant_vol_surface" that is
# inserted *before* it. The new sheet keeps the Tenors column (A) and the
# formulas/formatting exactly as-is, but every Quote in column B is
# replaced with a flat constant (10), and the unused helper column C is
# dropped.

$wb = $excel.ActiveWorkbook

# Grab the original sheet (it is currently the only / first sheet).
$volSheet = $wb.Worksheets.Item(1)

# Duplicate it in place, right before itself -- this clones data, number
# formats, column widths, borders, page setup, everything -- so the new
# sheet starts out pixel-identical to vol_surface.
$volSheet.Copy($volSheet)

# The freshly inserted copy is now the first sheet; rename it.
$constSheet = $wb.Worksheets.Item(1)
$constSheet.Name = "constant_vol_surface"

# Re-fetch vol_surface by name since indices shifted after the insert.
$volSheet = $wb.Worksheets.Item("vol_surface")

# Make the quotes column a flat constant of 10 on the new sheet.
for ($r = 2; $r -le 13; $r++) {
    $constSheet.Cells.Item($r, 2).Value2 = 10
}

# Drop the unused helper column (C) that only carried empty formatted
# cells on vol_surface.
$constSheet.Columns.Item(3).Delete()

# Match the selections / active sheet seen in the edited workbook:
# vol_surface is left with A1:B13 selected (no longer the active tab),
# constant_vol_surface becomes the active tab with B3:B13 selected.
[void]$volSheet.Range("A1:B13").Select()
[void]$constSheet.Range("B3:B13").Select()
$constSheet.Activate()
